$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.4916854829779206
$ws.Range("D2").Value = 51.53814122518933
$ws.Range("E2").Value = 0.1605476026824219
$ws.Range("C3").Value = 0.5903366173651952
$ws.Range("D3").Value = 47.04477903809251
$ws.Range("E3").Value = 0.2840161113617594
$ws.Range("C4").Value = 0.6482847852647085
$ws.Range("D4").Value = 47.05967344784042
$ws.Range("E4").Value = 0.5508174243493321
$ws.Range("C5").Value = 0.8139306930788575
$ws.Range("D5").Value = 42.7370915565283
$ws.Range("E5").Value = 1.254679736710638
$ws.Range("C6").Value = 1.221337463668168
$ws.Range("D6").Value = 38.64314332523863
$ws.Range("E6").Value = 2.187091025177575
$ws.Range("C7").Value = 1.560251443237467
$ws.Range("D7").Value = 33.47021790608454
$ws.Range("E7").Value = 3.106015921512773
$ws.Range("C8").Value = 2.039654110771267
$ws.Range("D8").Value = 25.64101013820287
$ws.Range("E8").Value = 4.333384461336562
$ws.Range("C9").Value = 2.559225855075067
$ws.Range("D9").Value = 20.2378741059077
$ws.Range("E9").Value = 5.017146462733881
$ws.Range("C10").Value = 2.754254870925581
$ws.Range("D10").Value = 15.09232158017083
$ws.Range("E10").Value = 5.989044024713238
$ws.Range("C11").Value = 3.144434948757379
$ws.Range("D11").Value = 10.00602065502239
$ws.Range("E11").Value = 6.722171313932866
$ws.Range("C12").Value = 3.453759085923031
$ws.Range("D12").Value = 7.036470305555672
$ws.Range("E12").Value = 7.295613231912729
$ws.Range("C13").Value = 3.692560634296318
$ws.Range("D13").Value = 3.973489613551741
$ws.Range("E13").Value = 7.912098565888122
$ws.Range("C14").Value = 3.763622124570412
$ws.Range("D14").Value = 1.528318616566812
$ws.Range("E14").Value = 8.124231356006836
$ws.Range("C15").Value = 3.830188535621611
$ws.Range("D15").Value = 1.962596700803358
$ws.Range("E15").Value = 8.21075264473329
$ws.Range("C16").Value = 3.944362221644016
$ws.Range("D16").Value = 1.633457232362763
$ws.Range("E16").Value = 8.178777427084865
$ws.Range("C17").Value = 3.907832339297824
$ws.Range("D17").Value = 0.9069630039691361
$ws.Range("E17").Value = 8.395343799053986
$ws.Range("C18").Value = 3.869295595706196
$ws.Range("D18").Value = -0.1374134823605665
$ws.Range("E18").Value = 8.390456777634949
$ws.Range("C19").Value = 3.900514804710857
$ws.Range("D19").Value = 0.01904213128491833
$ws.Range("E19").Value = 8.675144473929143
$ws.Range("C20").Value = 3.849040113610635
$ws.Range("D20").Value = 0.4595681123539536
$ws.Range("E20").Value = 8.639002345687908
$ws.Range("C21").Value = 3.945565585989047
$ws.Range("D21").Value = 0.4493284304179614
$ws.Range("E21").Value = 8.332535844734924
